# Applies the "Penalty Reward System" (unfinished) edit described in the diff:
#  - On the "Weekly Quantity" sheet, two weekly entries (rows 3 and 4, with
#    dates 45361.99999999999 and 45368.99999999999) are removed entirely,
#    shifting all following rows up by two. This shrinks the used range
#    from A1:B19 down to A1:B17.
#  - On the "Monthly Trend" sheet, the "Requested quantity" value for the
#    month in row 3 is changed from 330 to 190 (a partial / unfinished
#    update to match the weekly-sheet edit).

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# Remove the two obsolete weekly rows (old row 3 and old row 4). Deleting
# row 3 twice removes what were originally rows 3 and 4, shifting
# everything below up by two rows.
$weekly.Rows.Item(3).Delete()
$weekly.Rows.Item(3).Delete()

# Update the monthly trend value that was (partially) adjusted to reflect
# the weekly change above.
$monthly.Range("B3").Value = 190
